$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 29 de Agosto de 2020 a las 01:16"

# Re-label rows whose country changed position (shared-string reorder)
$ws.Range("A47").Value = "Japon"
$ws.Range("A48").Value = "Polonia"
$ws.Range("A123").Value = "Tunez"
$ws.Range("A124").Value = "Tailandia"
$ws.Range("A174").Value = "Polinesia Francesa"
$ws.Range("A175").Value = "Papua Nueva Guinea"
$ws.Range("A176").Value = "San Martin (Parte Holandesa)"
$ws.Range("A177").Value = "Burundi"
$ws.Range("A178").Value = "Comoras"
$ws.Range("A200").Value = "Islas Virgenes Britanicas"
$ws.Range("A201").Value = "Guam"
$ws.Range("A202").Value = "Fiyi"
$ws.Range("A203").Value = "Timor Oriental"

# Update numeric case data per row (B:Casos totales, C:Nuevos casos, D:Casos activos, E:Recuperados, F:Casos criticos, G:Muertes hoy, H:Muertes)
$ws.Range("B4").Value = 6092969
$ws.Range("C4").Value = 46335
$ws.Range("D4").Value = 3370495
$ws.Range("E4").Value = 2536639
$ws.Range("G4").Value = 1039
$ws.Range("H4").Value = 185835
$ws.Range("D5").Value = 2976796
$ws.Range("E5").Value = 716238
$ws.Range("B10").Value = 590520
$ws.Range("C10").Value = 8498
$ws.Range("D10").Value = 429620
$ws.Range("E10").Value = 142133
$ws.Range("G10").Value = 299
$ws.Range("H10").Value = 18767
$ws.Range("B14").Value = 392009
$ws.Range("C14").Value = 11717
$ws.Range("E14").Value = 103573
$ws.Range("G14").Value = 221
$ws.Range("H14").Value = 8271
$ws.Range("B35").Value = 93390
$ws.Range("C35").Value = 426
$ws.Range("D35").Value = 65285
$ws.Range("E35").Value = 26457
$ws.Range("G35").Value = 18
$ws.Range("H35").Value = 1648
$ws.Range("B36").Value = 90624
$ws.Range("C36").Value = 642
$ws.Range("D36").Value = 64602
$ws.Range("E36").Value = 24056
$ws.Range("G36").Value = 18
$ws.Range("H36").Value = 1966
$ws.Range("B47").Value = 65573
$ws.Range("C47").Value = 905
$ws.Range("D47").Value = 53995
$ws.Range("E47").Value = 10340
$ws.Range("G47").Value = 12
$ws.Range("H47").Value = 1238
$ws.Range("B48").Value = 65480
$ws.Range("C48").Value = 791
$ws.Range("D48").Value = 44785
$ws.Range("E48").Value = 18677
$ws.Range("G48").Value = 8
$ws.Range("H48").Value = 2018
$ws.Range("B53").Value = 53477
$ws.Range("C53").Value = 160
$ws.Range("D53").Value = 41017
$ws.Range("E53").Value = 11449
$ws.Range("B54").Value = 51113
$ws.Range("C54").Value = 357
$ws.Range("D54").Value = 47760
$ws.Range("E54").Value = 3164
$ws.Range("B90").Value = 10582
$ws.Range("C90").Value = 40
$ws.Range("E90").Value = 970
$ws.Range("B104").Value = 7012
$ws.Range("C104").Value = 19
$ws.Range("D104").Value = 6407
$ws.Range("E104").Value = 447
$ws.Range("B112").Value = 4700
$ws.Range("C112").Value = 2
$ws.Range("D112").Value = 1784
$ws.Range("B122").Value = 3697
$ws.Range("C122").Value = 46
$ws.Range("D122").Value = 2055
$ws.Range("E122").Value = 1621
$ws.Range("B123").Value = 3461
$ws.Range("C123").Value = 138
$ws.Range("D123").Value = 1522
$ws.Range("E123").Value = 1865
$ws.Range("G123").Value = 1
$ws.Range("H123").Value = 74
$ws.Range("B124").Value = 3410
$ws.Range("C124").Value = 6
$ws.Range("D124").Value = 3237
$ws.Range("E124").Value = 115
$ws.Range("H124").Value = 58
$ws.Range("B135").Value = 2471
$ws.Range("C135").Value = 56
$ws.Range("D135").Value = 1028
$ws.Range("E135").Value = 1337
$ws.Range("G135").Value = 1
$ws.Range("H135").Value = 106
$ws.Range("B149").Value = 1556
$ws.Range("C149").Value = 5
$ws.Range("D149").Value = 1352
$ws.Range("E149").Value = 161
$ws.Range("B150").Value = 1554
$ws.Range("C150").Value = 78
$ws.Range("D150").Value = 628
$ws.Range("E150").Value = 907
$ws.Range("G150").Value = 4
$ws.Range("H150").Value = 19
$ws.Range("B159").Value = 1175
$ws.Range("C159").Value = 2
$ws.Range("D159").Value = 1085
$ws.Range("E159").Value = 21
$ws.Range("B165").Value = 895
$ws.Range("C165").Value = 1
$ws.Range("D165").Value = 848
$ws.Range("E165").Value = 32
$ws.Range("B174").Value = 482
$ws.Range("C174").Value = 67
$ws.Range("D174").Value = 202
$ws.Range("E174").Value = 280
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 0
$ws.Range("B175").Value = 453
$ws.Range("C175").Value = 29
$ws.Range("D175").Value = 232
$ws.Range("E175").Value = 216
$ws.Range("G175").Value = 1
$ws.Range("H175").Value = 5
$ws.Range("B176").Value = 444
$ws.Range("C176").Value = 2
$ws.Range("D176").Value = 179
$ws.Range("E176").Value = 248
$ws.Range("H176").Value = 17
$ws.Range("B177").Value = 431
$ws.Range("D177").Value = 345
$ws.Range("E177").Value = 85
$ws.Range("H177").Value = 1
$ws.Range("B178").Value = 422
$ws.Range("D178").Value = 399
$ws.Range("E178").Value = 16
$ws.Range("H178").Value = 7
$ws.Range("B200").Value = 35
$ws.Range("C200").Value = 9
$ws.Range("D200").Value = 8
$ws.Range("E200").Value = 26
$ws.Range("B201").Value = 32
$ws.Range("D201").Value = 0
$ws.Range("E201").Value = 31
$ws.Range("H201").Value = 1
$ws.Range("B202").Value = 28
$ws.Range("D202").Value = 23
$ws.Range("E202").Value = 3
$ws.Range("H202").Value = 2
$ws.Range("B203").Value = 27
$ws.Range("D203").Value = 25
$ws.Range("E203").Value = 2
$ws.Range("H203").Value = 0
